$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Project name -> QuickMath, Sprint label -> Sprint 1
$ws.Range("A2").Value = "Sprint 1"
$ws.Range("A1").Value = "QuickMath"

# Remove the placeholder "User story " / "Work item " rows (4 and 5)
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()

# Update selection to G6
$ws.Range("G6").Select()
